$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @('frozenset({''IBM DB2''})', 'frozenset({''MS SQL Server''})', 0.04310344827586207, 0.03879310344827586, 0.02155172413793104, 0.5, 12.88888888888889, 0.01987960760998811, 1.922413793103448, 1, 1)
    3 = @('frozenset({''MS SQL Server''})', 'frozenset({''IBM DB2''})', 0.03879310344827586, 0.04310344827586207, 0.02155172413793104, 0.5555555555555556, 12.88888888888889, 0.01987960760998811, 2.15301724137931, 1, 1)
    12 = @('frozenset({''MS SQL Server''})', 'frozenset({''Oracle''})', 0.03879310344827586, 0.08620689655172414, 0.03017241379310345, 0.7777777777777778, 9.022222222222222, 0.0268281807372176, 4.112068965517241, 1, 1)
    13 = @('frozenset({''Oracle''})', 'frozenset({''MS SQL Server''})', 0.08620689655172414, 0.03879310344827586, 0.03017241379310345, 0.35, 9.022222222222222, 0.0268281807372176, 1.478779840848806, 1, 1)
    16 = @('frozenset({''IBM DB2''})', 'frozenset({''PostgreSQL''})', 0.04310344827586207, 0.09913793103448276, 0.03448275862068965, 0.7999999999999999, 8.069565217391304, 0.03020957193816885, 4.504310344827585, 1, 1)
    17 = @('frozenset({''PostgreSQL''})', 'frozenset({''IBM DB2''})', 0.09913793103448276, 0.04310344827586207, 0.03448275862068965, 0.3478260869565217, 8.069565217391304, 0.03020957193816885, 1.467241379310345, 1, 1)
    18 = @('frozenset({''Informix''})', 'frozenset({''HyperSQL''})', 0.02155172413793104, 0.1379310344827586, 0.02155172413793104, 1, 7.25, 0.01857907253269917, 'inf', 1, 1)
    19 = @('frozenset({''IBM DB2''})', 'frozenset({''HyperSQL''})', 0.04310344827586207, 0.1379310344827586, 0.04310344827586207, 1, 7.25, 0.03715814506539834, 'inf', 1, 1)
    40 = @('frozenset({''H2''})', 'frozenset({''IBM DB2''})', 0.1163793103448276, 0.04310344827586207, 0.02155172413793104, 0.1851851851851852, 4.296296296296297, 0.01653537455410226, 1.174373040752351, 1, 1)
    41 = @('frozenset({''MS Access''})', 'frozenset({''MySQL''})', 0.03017241379310345, 0.2327586206896552, 0.03017241379310345, 1, 4.296296296296297, 0.02314952437574316, 'inf', 1, 1)
    42 = @('frozenset({''MySQL''})', 'frozenset({''IBM DB2''})', 0.2327586206896552, 0.04310344827586207, 0.04310344827586207, 0.1851851851851852, 4.296296296296297, 0.03307074910820453, 1.174373040752351, 1, 1)
    43 = @('frozenset({''MySQL''})', 'frozenset({''MS Access''})', 0.2327586206896552, 0.03017241379310345, 0.03017241379310345, 0.1296296296296297, 4.296296296296297, 0.02314952437574316, 1.114269992663243, 1, 1)
    44 = @('frozenset({''IBM DB2''})', 'frozenset({''MySQL''})', 0.04310344827586207, 0.2327586206896552, 0.04310344827586207, 1, 4.296296296296297, 0.03307074910820453, 'inf', 1, 1)
    45 = @('frozenset({''MS SQL Server''})', 'frozenset({''MySQL''})', 0.03879310344827586, 0.2327586206896552, 0.03879310344827586, 1, 4.296296296296297, 0.02976367419738407, 'inf', 1, 1)
    46 = @('frozenset({''MySQL''})', 'frozenset({''MS SQL Server''})', 0.2327586206896552, 0.03879310344827586, 0.03879310344827586, 0.1666666666666667, 4.296296296296297, 0.02976367419738407, 1.153448275862069, 1, 1)
    47 = @('frozenset({''IBM DB2''})', 'frozenset({''H2''})', 0.04310344827586207, 0.1163793103448276, 0.02155172413793104, 0.5, 4.296296296296297, 0.01653537455410226, 1.767241379310345, 1, 1)
    52 = @('frozenset({''H2''})', 'frozenset({''PostgreSQL''})', 0.1163793103448276, 0.09913793103448276, 0.04310344827586207, 0.3703703703703704, 3.735909822866345, 0.03156584423305589, 1.43078093306288, 1, 1)
    53 = @('frozenset({''PostgreSQL''})', 'frozenset({''H2''})', 0.09913793103448276, 0.1163793103448276, 0.04310344827586207, 0.4347826086956522, 3.735909822866345, 0.03156584423305589, 1.563328912466843, 1, 1)
    58 = @('frozenset({''Oracle''})', 'frozenset({''MySQL''})', 0.08620689655172414, 0.2327586206896552, 0.06465517241379311, 0.75, 3.222222222222222, 0.04458977407847801, 3.068965517241379, 1, 1)
    60 = @('frozenset({''HyperSQL''})', 'frozenset({''H2''})', 0.1379310344827586, 0.1163793103448276, 0.05172413793103448, 0.375, 3.222222222222222, 0.0356718192627824, 1.413793103448276, 1, 1)
    61 = @('frozenset({''H2''})', 'frozenset({''HyperSQL''})', 0.1163793103448276, 0.1379310344827586, 0.05172413793103448, 0.4444444444444445, 3.222222222222222, 0.0356718192627824, 1.551724137931034, 1, 1)
}

foreach ($rowNum in $rowData.Keys) {
    $vals = $rowData[$rowNum]
    for ($col = 0; $col -lt $vals.Length; $col++) {
        $ws.Cells.Item([int]$rowNum, $col + 1).Value = $vals[$col]
    }
}
